$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 17.586354
$ws.Range("H2").Value = 52.759062
$ws.Range("I2").Value = 0.2178245326054132
$ws.Range("J2").Value = 0.2178245326054132
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.460132666666667
$ws.Range("N2").Value = 10.380398
$ws.Range("O2").Value = 0.01616897968344663
$ws.Range("P2").Value = 0.01616897968344663
$ws.Range("Q2").Value = 60.851117962964
$ws.Range("R2").Value = 547.660061666676
$ws.Range("S2").Value = 0.003522000442253184
$ws.Range("T2").Value = 0.003522000442253185

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 17.586354
$ws.Range("H3").Value = 52.759062
$ws.Range("I3").Value = 0.2178245326054132
$ws.Range("J3").Value = 0.2178245326054132
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 183.09167
$ws.Range("N3").Value = 549.27501
$ws.Range("O3").Value = 0.8555757185143522
$ws.Range("P3").Value = 0.8555757185143523
$ws.Range("Q3").Value = 3219.91492307118
$ws.Range("R3").Value = 28979.23430764062
$ws.Range("S3").Value = 0.1863653809939293
$ws.Range("T3").Value = 0.1863653809939294

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 17.586354
$ws.Range("H4").Value = 52.759062
$ws.Range("I4").Value = 0.2178245326054132
$ws.Range("J4").Value = 0.2178245326054132
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 27.44640466666667
$ws.Range("N4").Value = 82.339214
$ws.Range("O4").Value = 0.1282553018022011
$ws.Range("P4").Value = 0.1282553018022011
$ws.Range("Q4").Value = 482.682188495252
$ws.Range("R4").Value = 4344.139696457268
$ws.Range("S4").Value = 0.02793715116923065
$ws.Range("T4").Value = 0.02793715116923067

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 52.27042766666667
$ws.Range("H5").Value = 156.811283
$ws.Range("I5").Value = 0.6474213742983183
$ws.Range("J5").Value = 0.6474213742983183
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.460132666666667
$ws.Range("N5").Value = 10.380398
$ws.Range("O5").Value = 0.01616897968344663
$ws.Range("P5").Value = 0.01616897968344663
$ws.Range("Q5").Value = 180.8626142700705
$ws.Range("R5").Value = 1627.763528430634
$ws.Range("S5").Value = 0.0104681430476586
$ws.Range("T5").Value = 0.01046814304765861

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 52.27042766666667
$ws.Range("H6").Value = 156.811283
$ws.Range("I6").Value = 0.6474213742983183
$ws.Range("J6").Value = 0.6474213742983183
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 183.09167
$ws.Range("N6").Value = 549.27501
$ws.Range("O6").Value = 0.8555757185143522
$ws.Range("P6").Value = 0.8555757185143523
$ws.Range("Q6").Value = 9570.279893104203
$ws.Range("R6").Value = 86132.51903793782
$ws.Range("S6").Value = 0.553918007496833
$ws.Range("T6").Value = 0.5539180074968331

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 52.27042766666667
$ws.Range("H7").Value = 156.811283
$ws.Range("I7").Value = 0.6474213742983183
$ws.Range("J7").Value = 0.6474213742983183
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 27.44640466666667
$ws.Range("N7").Value = 82.339214
$ws.Range("O7").Value = 0.1282553018022011
$ws.Range("P7").Value = 0.1282553018022011
$ws.Range("Q7").Value = 1434.635309839063
$ws.Range("R7").Value = 12911.71778855156
$ws.Range("S7").Value = 0.08303522375382659
$ws.Range("T7").Value = 0.0830352237538266

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.87955133333333
$ws.Range("H8").Value = 32.638654
$ws.Range("I8").Value = 0.1347540930962685
$ws.Range("J8").Value = 0.1347540930962685
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.460132666666667
$ws.Range("N8").Value = 10.380398
$ws.Range("O8").Value = 0.01616897968344663
$ws.Range("P8").Value = 0.01616897968344663
$ws.Range("Q8").Value = 37.64469096714355
$ws.Range("R8").Value = 338.802218704292
$ws.Range("S8").Value = 0.00217883619353484
$ws.Range("T8").Value = 0.002178836193534841

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.87955133333333
$ws.Range("H9").Value = 32.638654
$ws.Range("I9").Value = 0.1347540930962685
$ws.Range("J9").Value = 0.1347540930962685
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 183.09167
$ws.Range("N9").Value = 549.27501
$ws.Range("O9").Value = 0.8555757185143522
$ws.Range("P9").Value = 0.8555757185143523
$ws.Range("Q9").Value = 1991.955222470727
$ws.Range("R9").Value = 17927.59700223654
$ws.Range("S9").Value = 0.1152923300235898
$ws.Range("T9").Value = 0.1152923300235898

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.87955133333333
$ws.Range("H10").Value = 32.638654
$ws.Range("I10").Value = 0.1347540930962685
$ws.Range("J10").Value = 0.1347540930962685
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 27.44640466666667
$ws.Range("N10").Value = 82.339214
$ws.Range("O10").Value = 0.1282553018022011
$ws.Range("P10").Value = 0.1282553018022011
$ws.Range("Q10").Value = 298.6045684864396
$ws.Range("R10").Value = 2687.441116377956
$ws.Range("S10").Value = 0.01728292687914381
$ws.Range("T10").Value = 0.01728292687914382
